# Auto-generated Excel COM-interop script applying the cryptos list update
# (values refreshed per commit "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.853.29'
$ws.Range("E2").Value = '  +3.35%  '

$ws.Range("D3").Value = '3.413.92'
$ws.Range("E3").Value = '  +3.57%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.44'
$ws.Range("E5").Value = '  +3.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.67'
$ws.Range("E6").Value = '  +8.10%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.413.69'
$ws.Range("E8").Value = '  +3.60%  '

$ws.Range("E9").Value = '  +1.82%  '

$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.49'
$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.128'
$ws.Range("E11").Value = '  +10.19%  '

$ws.Range("E12").Value = '  +7.11%  '

$ws.Range("D13").Value = '3.998.76'
$ws.Range("E13").Value = '  +3.52%  '

$ws.Range("E14").Value = '  +2.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +8.78%  '

$ws.Range("D16").Value = '3.410.18'
$ws.Range("E16").Value = '  +3.36%  '

$ws.Range("E17").Value = '  +6.53%  '

$ws.Range("D18").Value = '61.944.33'
$ws.Range("E18").Value = '  +3.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.13'
$ws.Range("E19").Value = '  +7.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.90'
$ws.Range("E20").Value = '  +5.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.52'
$ws.Range("E21").Value = '  +8.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.92'
$ws.Range("E22").Value = '  +12.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.573'
$ws.Range("E23").Value = '  +4.23%  '

$ws.Range("D24").Value = '3.551.53'
$ws.Range("E24").Value = '  +3.65%  '

$ws.Range("E25").Value = '  +19.69%  '

$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.67'
$ws.Range("E27").Value = '  +4.80%  '

$ws.Range("E28").Value = '  +11.05%  '

$ws.Range("E29").Value = '  +5.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.32'
$ws.Range("E31").Value = '  +6.98%  '

$ws.Range("E32").Value = '  +6.02%  '

$ws.Range("E33").Value = '  +4.40%  '

$ws.Range("D34").Value = '3.445.72'
$ws.Range("E34").Value = '  +3.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.60'
$ws.Range("E36").Value = '  +4.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.49'
$ws.Range("E37").Value = '  +4.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.00'
$ws.Range("E38").Value = '  +4.18%  '

$ws.Range("E39").Value = '  +6.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.74'
$ws.Range("E40").Value = '  +3.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0792'

$ws.Range("E42").Value = '  +15.51%  '

$ws.Range("E43").Value = '  +7.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.32'
$ws.Range("E44").Value = '  +12.67%  '

$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("E46").Value = '  +6.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.48'
$ws.Range("E47").Value = '  +4.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.65'
$ws.Range("E48").Value = '  +3.38%  '

$ws.Range("E49").Value = '  +4.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.20'
$ws.Range("E50").Value = '  +6.66%  '

$ws.Range("D51").Value = '2.393.16'
$ws.Range("E51").Value = '  +11.19%  '
